# Auto-applies the cryptos.xlsx price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$CellRef, [string]$Text)
    $cell = $Worksheet.Range($CellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = '@'
    $cell.Value = $Text
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws "D2" '27.712.12'
Set-TextValue $ws "E2" '  -0.31%  '

# Row 3
Set-TextValue $ws "D3" '1.847.04'
Set-TextValue $ws "E3" '  -0.84%  '

# Row 4
Set-TextValue $ws "D4" '1.013'
Set-TextValue $ws "E4" '  -2.01%  '

# Row 5
Set-TextValue $ws "D5" '319.28'
Set-TextValue $ws "E5" '  -0.97%  '

# Row 6
Set-TextValue $ws "D6" '1.011'
Set-TextValue $ws "E6" '  -1.81%  '

# Row 7
Set-TextValue $ws "D7" '0.4314'
Set-TextValue $ws "E7" '  -2.40%  '

# Row 8
Set-TextValue $ws "D8" '0.3750'
Set-TextValue $ws "E8" '  -1.55%  '

# Row 9
Set-TextValue $ws "D9" '0.07352'
Set-TextValue $ws "E9" '  -1.07%  '

# Row 10
Set-TextValue $ws "D10" '0.8775'
Set-TextValue $ws "E10" '  -1.02%  '

# Row 11
Set-TextValue $ws "D11" '21.63'
Set-TextValue $ws "E11" '  -0.68%  '

# Row 12
Set-TextValue $ws "D12" '1.857.21'
Set-TextValue $ws "E12" '  -0.58%  '

# Row 13
Set-TextValue $ws "D13" '6.725'
Set-TextValue $ws "E13" '  -0.65%  '

# Row 14
Set-TextValue $ws "D14" '5.446'
Set-TextValue $ws "E14" '  -1.87%  '

# Row 15
Set-TextValue $ws "D15" '0.07132'
Set-TextValue $ws "E15" '  -0.72%  '

# Row 16
Set-TextValue $ws "D16" '89.14'
Set-TextValue $ws "E16" '  +5.61%  '

# Row 17
Set-TextValue $ws "D17" '1.015'
Set-TextValue $ws "E17" '  -1.93%  '

# Row 18
Set-TextValue $ws "D18" '0.000008997'
Set-TextValue $ws "E18" '  -1.33%  '

# Row 19
Set-TextValue $ws "D19" '1.009'
Set-TextValue $ws "E19" '  -2.02%  '

# Row 20
Set-TextValue $ws "D20" '15.48'
Set-TextValue $ws "E20" '  -0.04%  '

# Row 21
Set-TextValue $ws "D21" '27.697.01'
Set-TextValue $ws "E21" '  -0.44%  '

# Row 22
Set-TextValue $ws "D22" '5.219'
Set-TextValue $ws "E22" '  -1.70%  '

# Row 23
Set-TextValue $ws "E23" '  -1.89%  '

# Row 24
Set-TextValue $ws "D24" '2.073.51'
Set-TextValue $ws "E24" '  -0.92%  '

# Row 25
Set-TextValue $ws "D25" '2.007'
Set-TextValue $ws "E25" '  -0.67%  '

# Row 26
Set-TextValue $ws "D26" '155.28'
Set-TextValue $ws "E26" '  -1.86%  '

# Row 27
Set-TextValue $ws "D27" '18.66'
Set-TextValue $ws "E27" '  -1.26%  '

# Row 28
Set-TextValue $ws "D28" '2.198'
Set-TextValue $ws "E28" '  +10.86%  '

# Row 29
Set-TextValue $ws "D29" '5.389'
Set-TextValue $ws "E29" '  +0.11%  '

# Row 30
Set-TextValue $ws "D30" '119.39'
Set-TextValue $ws "E30" '  +0.52%  '

# Row 31
Set-TextValue $ws "D31" '0.08943'
Set-TextValue $ws "E31" '  -0.97%  '

# Row 32
Set-TextValue $ws "D32" '1.232'
Set-TextValue $ws "E32" '  -0.22%  '

# Row 33
Set-TextValue $ws "D33" '0.7795'
Set-TextValue $ws "E33" '  +0.20%  '

# Row 34
Set-TextValue $ws "D34" '4.555'
Set-TextValue $ws "E34" '  -0.89%  '

# Row 35
Set-TextValue $ws "D35" '2.893'
Set-TextValue $ws "E35" '  -3.38%  '

# Row 36
Set-TextValue $ws "E36" '  -1.91%  '

# Row 37
Set-TextValue $ws "D37" '1.134'
Set-TextValue $ws "E37" '  -0.81%  '

# Row 38
Set-TextValue $ws "D38" '0.05344'

# Row 39
Set-TextValue $ws "E39" '  -1.04%  '

# Row 40
Set-TextValue $ws "D40" '7.337'
Set-TextValue $ws "E40" '  +6.29%  '

# Row 41
Set-TextValue $ws "D41" '2.890'
Set-TextValue $ws "E41" '  +0.36%  '

# Row 42
Set-TextValue $ws "B42" 'Algorand'
Set-TextValue $ws "C42" 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws "D42" '0.1687'
Set-TextValue $ws "E42" '  -0.37%  '

# Row 43
Set-TextValue $ws "B43" 'TheSandbox'
Set-TextValue $ws "C43" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws "D43" '0.5121'
Set-TextValue $ws "E43" '  -1.73%  '

# Row 44
Set-TextValue $ws "D44" '8.807'
Set-TextValue $ws "E44" '  +0.93%  '

# Row 45
Set-TextValue $ws "B45" 'EnergySwap'
Set-TextValue $ws "C45" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws "D45" '10.72'
Set-TextValue $ws "E45" '  -0.01%  '

# Row 46
Set-TextValue $ws "B46" 'Quant'
Set-TextValue $ws "C46" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws "D46" '109.14'
Set-TextValue $ws "E46" '  -2.07%  '

# Row 47
Set-TextValue $ws "D47" '0.4759'
Set-TextValue $ws "E47" '  +0.63%  '

# Row 48
Set-TextValue $ws "D48" '0.06475'
Set-TextValue $ws "E48" '  -2.71%  '

# Row 49
Set-TextValue $ws "B49" 'NEARProtocol'
Set-TextValue $ws "C49" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws "D49" '1.692'
Set-TextValue $ws "E49" '  -1.18%  '

# Row 50
Set-TextValue $ws "B50" 'PaxDollar'
Set-TextValue $ws "C50" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws "D50" '1.012'
Set-TextValue $ws "E50" '  -1.82%  '

# Row 51
Set-TextValue $ws "D51" '1.850'
Set-TextValue $ws "E51" '  -3.86%  '
